$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (Strikeouts) values regenerated for column G, rows 2-15
$gValues = @{
    2  = 5
    3  = 6
    4  = 2
    5  = 3
    6  = 4
    7  = 2
    8  = 7
    9  = 3
    10 = 7
    11 = 3
    12 = 4
    13 = 2
    14 = 1
    15 = 2
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
